{"js": "// Translate the English source strings in the Word document to their\n// Afrikaans equivalents, one search-and-replace per paragraph run.\nconst replacements = [\n  [\n    \"Appendix 15: SWIFT Interview Guide: Referral Services \",\n    \"Bylaag 15: SWIFT Onderhoudgids: Verwysingsdienste \",\n  ],\n  [\"Briefing:\", \"Inligtingsessie:\"],\n  [\"What did you find helpful? \", \"Wat het jy nuttig gevind? \"],\n  [\n    \"Was there anything more you wish it had? \",\n    \"Was daar enige iets meer wat jy wens dit het gehad? \",\n  ],\n  [\n    \"Did you find out about any services which you didn\\u2019t know about before?\",\n    \"Het jy uitgevind van dienste wat jy nie voorheen van geweet het nie?\",\n  ],\n  [\n    \"Probe: What could have made them more helpful?\",\n    \"Ondersoek: Wat sou hulle meer nuttig gemaak het?\",\n  ],\n  [\n    \"Do you have any suggestions for how we could improve the \\u2018help\\u2019 menu in the programme?\",\n    \"Het jy enige voorstelle oor hoe ons die help-kieslys kan verbeter in die program?\",\n  ],\n  [\"Debriefing\", \"Afsluiting\"],\n  [\n    \"Is there anything that we haven\\u2019t spoken about that you\\u2019d like us to know? \",\n    \"Is daar enige iets wat ons nog nie bespreek het nie wat jy graag wil h\\u00ea ons moet weet? \",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the English source strings in the Word document to their\n# Afrikaans equivalents, one find-and-replace per paragraph run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"Appendix 15: SWIFT Interview Guide: Referral Services \", \"Bylaag 15: SWIFT Onderhoudgids: Verwysingsdienste \"),\n    @(\"Briefing:\", \"Inligtingsessie:\"),\n    @(\"What did you find helpful? \", \"Wat het jy nuttig gevind? \"),\n    @(\"Was there anything more you wish it had? \", \"Was daar enige iets meer wat jy wens dit het gehad? \"),\n    @(\"Did you find out about any services which you didn\u2019t know about before?\", \"Het jy uitgevind van dienste wat jy nie voorheen van geweet het nie?\"),\n    @(\"Probe: What could have made them more helpful?\", \"Ondersoek: Wat sou hulle meer nuttig gemaak het?\"),\n    @(\"Do you have any suggestions for how we could improve the \u2018help\u2019 menu in the programme?\", \"Het jy enige voorstelle oor hoe ons die help-kieslys kan verbeter in die program?\"),\n    @(\"Debriefing\", \"Afsluiting\"),\n    @(\"Is there anything that we haven\u2019t spoken about that you\u2019d like us to know? \", \"Is daar enige iets wat ons nog nie bespreek het nie wat jy graag wil h\u00ea ons moet weet? \")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
